$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 9813.61227414258
$ws.Range("F2").Value = -46.2024064708421

$ws.Range("C3").Value = 9757.80293385522
$ws.Range("F3").Value = 283.227836375068

$ws.Range("C4").Value = 9367.06869931902
$ws.Range("F4").Value = 267.185440346171

$ws.Range("C5").Value = 7269.48020724183
$ws.Range("F5").Value = 176.68169116291

$ws.Range("C6").Value = 7658.01890033183
$ws.Range("F6").Value = 212.605437220458

$ws.Range("C7").Value = 10166.6111885471
$ws.Range("F7").Value = 315.518622054275
